$d = $word.ActiveDocument

# Locate the paragraph that ends the "cadastro de contatos" requirement -
# "Tipo de Contato que podem ser Confirmacao de recebimento, ... data
# posterior ao contato efetuado." (the paragraph that also carries the
# _GoBack bookmark). It is the last paragraph of the document body.
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -like "*data posterior ao contato efetuado*") {
        $targetPara = $candidate
    }
}

# Place the insertion point at the very end of that paragraph (after the
# trailing space and the bookmark, right before the paragraph mark) and
# press Enter to start a brand new paragraph below it, exactly like a user
# finishing the existing requirement text and continuing with a new one.
$insertionPoint = $targetPara.Range
$insertionPoint.Collapse(0)
$insertionPoint.InsertParagraphAfter()

# The freshly created paragraph is now the new last paragraph of the body.
$newPara = $d.Paragraphs.Last
$newRange = $newPara.Range
$newRange.Collapse(0)
$newRange.InsertAfter("O sistema deve ter um controle de movimento de caixa, para registro de recebimentos que podem acontecer em loja. Também deve ser possível acompanhar as contas a receber da empresa de cobrança, assim como suas contas a pagar, e emitir um relatório gerencial de fluxo de caixa.")

# Split the new paragraph's text into its two sentences as two separate
# runs (mirroring how the text was authored over two edits).
$fullRange = $newPara.Range

$firstSentence = $fullRange.Duplicate
$firstSentence.Find.Execute("O sistema deve ter um controle de movimento de caixa, para registro de recebimentos que podem acontecer em loja. ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$firstSentence.Bold = 1
$firstSentence.Bold = 0

$secondSentence = $fullRange.Duplicate
$secondSentence.Find.Execute("Também deve ser possível acompanhar as contas a receber da empresa de cobrança, assim como suas contas a pagar, e emitir um relatório gerencial de fluxo de caixa.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$secondSentence.Bold = 1
$secondSentence.Bold = 0

# The section's page is explicitly stamped as portrait orientation.
$d.PageSetup.Orientation = 0
